$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 429 - this shifts rows 429..531 down to 430..532,
# matching the dimension growing from A1:R531 to A1:R532.
$ws.Rows.Item(429).Insert()

# Fill in the new row 429 with the new data record.
$ws.Cells.Item(429, 1).Value = 3
$ws.Cells.Item(429, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(429, 3).Value = "Coquimbo"
$ws.Cells.Item(429, 4).Value = [datetime]"2023-10-12"
$ws.Cells.Item(429, 5).Value = 5
$ws.Cells.Item(429, 6).Value = 100112001
$ws.Cells.Item(429, 7).Value = "Berenjena"
$ws.Cells.Item(429, 8).Value = "Sin especificar"
$ws.Cells.Item(429, 9).Value = "Primera"
$ws.Cells.Item(429, 10).Value = 100
$ws.Cells.Item(429, 11).Value = 8500
$ws.Cells.Item(429, 12).Value = 9000
$ws.Cells.Item(429, 13).Value = 8750
$ws.Cells.Item(429, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(429, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(429, 16).Value = 146
$ws.Cells.Item(429, 17).Value = 60
$ws.Cells.Item(429, 18).Value = "Hortaliza"
